$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to the new columns C and D before writing values,
# so strings like "010" / "01" aren't coerced into numbers.
$ws.Range("C1:C22").NumberFormat = "@"
$ws.Range("D1:D22").NumberFormat = "@"

# Column C: "010", "020", ... "220" (row 1..22), written first so the
# shared-string table picks these up before column D's values.
for ($i = 1; $i -le 22; $i++) {
    $code = [string]($i * 10)
    $code = $code.PadLeft(3, '0')
    $ws.Cells.Item($i, 3).Value = $code
}

# Column D: "01", "02", ... "22" (row 1..22), written after column C.
for ($i = 1; $i -le 22; $i++) {
    $code2 = [string]$i
    $code2 = $code2.PadLeft(2, '0')
    $ws.Cells.Item($i, 4).Value = $code2
}

# Match the authored column widths for the two new columns.
$ws.Columns.Item(3).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(4).ColumnWidth = 2.1666666666666665

# Move the active selection as recorded in the saved workbook.
$ws.Range("G18").Select()
